$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(46, 1).Value = "2025-04-29 04:58:08"
$ws.Cells.Item(46, 2).Value = 148
